$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared-string value in E2: "AEC4" -> "AEC1"
$ws.Range("E2").Value = "AEC1"

# Update numeric values in row 1
$ws.Range("D1").Value = 100
$ws.Range("E1").Value = 125

# Move the active selection to E1 (as in the final sheetView)
$ws.Activate()
$ws.Range("E1").Select()

$wb.Save()
